# Adds the new weekly backlog rows (semana 5, 03/02/2025-07/02/2025) that were
# appended to the "SPN" and "ITI" sheets of Base/consolidado.xlsx.
#
# Columns: A Setor | B Responsavel | C Ano | D Semana | E Inicio_Semana |
#          F Final_Semana | G Incidente | H Backlog | I Data | J Status | K Coordenador
#
# E, F, H and I hold day-first date-like text (e.g. "03/02/2025", "01/2025")
# that must stay plain text, not get auto-converted to Excel date serials.
# A leading apostrophe forces the text interpretation; the style is then put
# back to "Normal" so no extra number-format style gets attached to the cell
# (matching the source file, where these new cells carry no explicit style).
#
# NOTE: worksheet COM objects are looped over inline (not passed into a
# function) because passing a Range/Worksheet reference as a function
# parameter in this host loses the underlying COM binding.

$wb = $excel.ActiveWorkbook

$spnRows = @(
    @(16, 'SPN', 'Higor Cruz', 2025, 5, '03/02/2025', '07/02/2025', 322731, '01/2025', '03/02/2025', 'Pendente', 'Willian Jones'),
    @(17, 'SPN', 'Higor Cruz', 2025, 5, '03/02/2025', '07/02/2025', 322526, '01/2025', '03/02/2025', 'Pendente', 'Willian Jones'),
    @(18, 'SPN', 'Higor Cruz', 2025, 5, '03/02/2025', '07/02/2025', 321751, '01/2025', '03/02/2025', 'Pendente', 'Willian Jones'),
    @(19, 'SPN', 'Higor Cruz', 2025, 5, '03/02/2025', '07/02/2025', 321760, '01/2025', '03/02/2025', 'Pendente', 'Willian Jones'),
    @(20, 'SPN', 'Arthur Hassuma', 2025, 5, '03/02/2025', '07/02/2025', 322877, '01/2025', '03/02/2025', 'Pendente', 'Willian Jones'),
    @(21, 'SPN', 'Luan Pierry', 2025, 5, '03/02/2025', '07/02/2025', 322586, '01/2025', '03/02/2025', 'Pendente', 'Willian Jones'),
    @(22, 'SPN', 'Luan Pierry', 2025, 5, '03/02/2025', '07/02/2025', 322053, '01/2025', '03/02/2025', 'Pendente', 'Willian Jones'),
    @(23, 'SPN', 'Mara Neves', 2025, 5, '03/02/2025', '07/02/2025', 322696, '01/2025', '03/02/2025', 'Pendente', 'Willian Jones'),
    @(24, 'SPN', 'Mara Neves', 2025, 5, '03/02/2025', '07/02/2025', 322164, '01/2025', '03/02/2025', 'Pendente', 'Willian Jones')
)

$itiRows = @(
    @(39, 'ITI', 'Lourival Moizés', 2025, 5, '03/02/2025', '07/02/2025', 322346, '01/2025', '03/02/2025', 'Pendente', 'Emerson Simette'),
    @(40, 'ITI', 'Guilherme Worel', 2025, 5, '03/02/2025', '07/02/2025', 321835, '01/2025', '03/02/2025', 'Pendente', 'Emerson Simette'),
    @(41, 'ITI', 'Guilherme Worel', 2025, 5, '03/02/2025', '07/02/2025', 322897, '01/2025', '03/02/2025', 'Pendente', 'Emerson Simette'),
    @(42, 'ITI', 'Guilherme Worel', 2025, 5, '03/02/2025', '07/02/2025', 322991, '01/2025', '03/02/2025', 'Pendente', 'Emerson Simette'),
    @(43, 'ITI', 'Jorgenaldo Reis', 2025, 5, '03/02/2025', '07/02/2025', 322655, '01/2025', '03/02/2025', 'Pendente', 'Emerson Simette'),
    @(44, 'ITI', 'Jose Acevedo', 2025, 5, '03/02/2025', '07/02/2025', 322167, '01/2025', '03/02/2025', 'Pendente', 'Emerson Simette'),
    @(45, 'ITI', 'Erick da Silva', 2025, 5, '03/02/2025', '07/02/2025', 322927, '01/2025', '03/02/2025', 'Pendente', 'Emerson Simette'),
    @(46, 'ITI', 'Erick da Silva', 2025, 5, '03/02/2025', '07/02/2025', 322759, '01/2025', '03/02/2025', 'Pendente', 'Emerson Simette'),
    @(47, 'ITI', 'Erick da Silva', 2025, 5, '03/02/2025', '07/02/2025', 322764, '01/2025', '03/02/2025', 'Pendente', 'Emerson Simette'),
    @(48, 'ITI', 'Erick da Silva', 2025, 5, '03/02/2025', '07/02/2025', 322804, '01/2025', '03/02/2025', 'Pendente', 'Emerson Simette'),
    @(49, 'ITI', 'Lourival Moizés', 2025, 5, '03/02/2025', '07/02/2025', 321811, '01/2025', '03/02/2025', 'Pendente', 'Emerson Simette')
)

$ws = $wb.Worksheets.Item("SPN")
$rows = $spnRows
foreach ($row in $rows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]

    $ws.Cells.Item($r, 5).Value = "'" + $row[5]
    $ws.Cells.Item($r, 5).Style = "Normal"

    $ws.Cells.Item($r, 6).Value = "'" + $row[6]
    $ws.Cells.Item($r, 6).Style = "Normal"

    $ws.Cells.Item($r, 7).Value = $row[7]

    $ws.Cells.Item($r, 8).Value = "'" + $row[8]
    $ws.Cells.Item($r, 8).Style = "Normal"

    $ws.Cells.Item($r, 9).Value = "'" + $row[9]
    $ws.Cells.Item($r, 9).Style = "Normal"

    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
}

$ws = $wb.Worksheets.Item("ITI")
$rows = $itiRows
foreach ($row in $rows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]

    $ws.Cells.Item($r, 5).Value = "'" + $row[5]
    $ws.Cells.Item($r, 5).Style = "Normal"

    $ws.Cells.Item($r, 6).Value = "'" + $row[6]
    $ws.Cells.Item($r, 6).Style = "Normal"

    $ws.Cells.Item($r, 7).Value = $row[7]

    $ws.Cells.Item($r, 8).Value = "'" + $row[8]
    $ws.Cells.Item($r, 8).Style = "Normal"

    $ws.Cells.Item($r, 9).Value = "'" + $row[9]
    $ws.Cells.Item($r, 9).Style = "Normal"

    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
}
